$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '241.82'
Set-TextCell 2 7 '19'
Set-TextCell 3 4 '21.84'
Set-TextCell 3 7 '19'
Set-TextCell 4 4 '5.375'
Set-TextCell 4 7 '19'
Set-TextCell 5 4 '0.05681'
Set-TextCell 5 7 '19'
Set-TextCell 6 4 '3.404'
Set-TextCell 6 7 '19'
Set-TextCell 7 4 '6.285'
Set-TextCell 7 7 '19'
Set-TextCell 8 2 'FTXToken'
Set-TextCell 8 3 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell 8 4 '1.137'
Set-TextCell 8 5 '7FTXTokenFTT'
Set-TextCell 8 7 '19'
Set-TextCell 9 2 'MXToken'
Set-TextCell 9 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 9 4 '0.8057'
Set-TextCell 9 5 '8MXTokenMX'
Set-TextCell 9 7 '19'
Set-TextCell 10 2 'WazirX'
Set-TextCell 10 3 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell 10 4 '0.1419'
Set-TextCell 10 5 '9WazirXWRX'
Set-TextCell 10 7 '19'
Set-TextCell 11 2 'MandalaExchangeToken'
Set-TextCell 11 3 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell 11 4 '0.07274'
Set-TextCell 11 5 '10MandalaExchangeTokenMDX'
Set-TextCell 11 7 '19'
Set-TextCell 12 2 'LiechtensteinCryptoassetsExchange'
Set-TextCell 12 3 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell 12 4 '0.03039'
Set-TextCell 12 5 '11LiechtensteinCryptoassetsExchangeLCX'
Set-TextCell 12 7 '19'
Set-TextCell 13 2 'BitrueCoin'
Set-TextCell 13 3 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell 13 4 '0.03113'
Set-TextCell 13 5 '12BitrueCoinBTR'
Set-TextCell 13 7 '19'
Set-TextCell 14 2 'BitMartToken'
Set-TextCell 14 3 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell 14 4 '0.09350'
Set-TextCell 14 5 '13BitMartTokenBMX'
Set-TextCell 14 7 '19'
Set-TextCell 15 2 'MCDex'
Set-TextCell 15 3 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextCell 15 4 '3.915'
Set-TextCell 15 5 '14MCDexMCB'
Set-TextCell 15 7 '19'
Set-TextCell 16 2 'BitForexToken'
Set-TextCell 16 3 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell 16 4 '0.001590'
Set-TextCell 16 5 '15BitForexTokenBF'
Set-TextCell 16 7 '19'
Set-TextCell 17 2 'CoinExToken'
Set-TextCell 17 3 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextCell 17 4 '0.04804'
Set-TextCell 17 5 '16CoinExTokenCET'
Set-TextCell 17 7 '19'
Set-TextCell 18 2 'One'
Set-TextCell 18 3 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextCell 18 4 '0.0005820'
Set-TextCell 18 5 '17OneONE'
Set-TextCell 18 7 '19'
Set-TextCell 19 2 'TigerCash'
Set-TextCell 19 3 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell 19 4 '0.006326'
Set-TextCell 19 5 '18TigerCashTCH'
Set-TextCell 19 7 '19'
Set-TextCell 20 2 'HotbitToken'
Set-TextCell 20 3 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextCell 20 4 '0.004057'
Set-TextCell 20 5 '19HotbitTokenHTB'
Set-TextCell 20 7 '19'
Set-TextCell 21 2 'BitKan'
Set-TextCell 21 3 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextCell 21 4 '0.0009938'
Set-TextCell 21 5 '20BitKanKAN'
Set-TextCell 21 7 '19'
Set-TextCell 22 2 'NitroEx'
Set-TextCell 22 3 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-TextCell 22 4 '0.0001500'
Set-TextCell 22 5 '21NitroExNTX'
Set-TextCell 22 7 '19'
Set-TextCell 23 2 'LEO'
Set-TextCell 23 3 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell 23 4 '3.735'
Set-TextCell 23 5 '22LEOLEO'
Set-TextCell 23 7 '19'
Set-TextCell 24 2 'BTSEToken'
Set-TextCell 24 3 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextCell 24 4 '2.150'
Set-TextCell 24 5 '23BTSETokenBTSE'
Set-TextCell 24 7 '19'
Set-TextCell 25 2 'BitpandaEcosystemToken'
Set-TextCell 25 3 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextCell 25 4 '0.3258'
Set-TextCell 25 5 '24BitpandaEcosystemTokenBEST'
Set-TextCell 25 7 '19'
Set-TextCell 26 2 'ProBitToken'
Set-TextCell 26 3 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextCell 26 4 '0.1310'
Set-TextCell 26 5 '25ProBitTokenPROB'
Set-TextCell 26 7 '19'
Set-TextCell 27 4 '0.0003999'
Set-TextCell 27 5 '26UpBotsUBXT'
Set-TextCell 27 7 '19'
Set-TextCell 28 7 '19'
Set-TextCell 29 7 '19'
Set-TextCell 30 7 '19'
Set-TextCell 31 7 '19'
Set-TextCell 32 7 '19'
Set-TextCell 33 7 '19'
Set-TextCell 34 7 '19'
Set-TextCell 35 7 '19'
Set-TextCell 36 7 '19'
Set-TextCell 37 7 '19'
Set-TextCell 38 7 '19'
Set-TextCell 39 7 '19'
Set-TextCell 40 4 '0.03808'
Set-TextCell 40 7 '19'
Set-TextCell 41 4 '0.006684'
Set-TextCell 41 5 '40KickTokenKICKBestin24h'
Set-TextCell 41 7 '19'
Set-TextCell 42 4 '0.1048'
Set-TextCell 42 7 '19'
Set-TextCell 43 4 '0.003000'
Set-TextCell 43 7 '19'
Set-TextCell 44 4 '0.006475'
Set-TextCell 44 7 '19'
Set-TextCell 45 4 '0.00005597'
Set-TextCell 45 7 '19'
Set-TextCell 46 7 '19'
Set-TextCell 47 4 '0.3900'
Set-TextCell 47 7 '19'
Set-TextCell 48 5 '47BOLOBOLOWorstin24h'
Set-TextCell 48 7 '19'
Set-TextCell 49 7 '19'
Set-TextCell 50 7 '19'
Set-TextCell 51 7 '19'
